$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the DataDir/^ROMCS(O) pair with the NC_A/NC_B pair (rows 38-39 <-> 40-41),
# and rename the NC_A/NC_B pins to RX/TX for the new ESP-01S header.
# Columns B (pin number) and D (edge connector) stay put; A (name), C (colour)
# and E (designator) move with the signal.
$ws.Range("A38").Value = "RX"
$ws.Range("C38").Value = "Green"
$ws.Range("E38").Value = "L2"

$ws.Range("A39").Value = "TX"
$ws.Range("C39").Value = "Blue"
$ws.Range("E39").Value = "U2"

$ws.Range("A40").Value = "DataDir"
$ws.Range("C40").Value = "Orange"
$ws.Range("E40").Value = "U1"

$ws.Range("A41").Value = "^ROMCS (O)"
$ws.Range("C41").Value = "Yellow"
$ws.Range("E41").Value = "L1"

# Update the active selection left behind in the saved view state.
$ws.Range("D53").Select()
